# Update Tgfb3-Eng.xlsx NATMI edge table with newly recomputed TPM-derived
# ligand/receptor expression values and all values that are derived from them
# (totals, specificities, and edge weights/specificities).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values per row (only cells whose value actually changes are listed).
$updates = @{
    2  = @{ G=1.103903333333333; H=3.31171;            I=0.02393122995918198; J=0.02393122995918198;
             M=135.0916853333333; N=405.2750559999999; O=0.7123704212620513;  P=0.7123704212620514;
             Q=149.1281617450844; R=1342.15345570576;  S=0.01704790036734149; T=0.01704790036734149 }
    3  = @{ G=1.103903333333333; H=3.31171;            I=0.02393122995918198; J=0.02393122995918198;
             O=0.2125756143240238; P=0.2125756143240238;
             Q=44.50073957283556; R=400.50665615552;   S=0.005087195910102591; T=0.005087195910102592 }
    4  = @{ G=1.103903333333333; H=3.31171;            I=0.02393122995918198; J=0.02393122995918198;
             M=14.23299766666667; N=42.698993;         O=0.07505396441392481; P=0.07505396441392483;
             Q=15.71185356755889; R=141.40668210803;   S=0.001796133681737896; T=0.001796133681737896 }
    5  = @{ H=70.73212899999999; I=0.5111277390235027; J=0.5111277390235027;
             M=135.0916853333333; N=405.2750559999999; O=0.7123704212620513;  P=0.7123704212620514;
             Q=3185.107504608246; R=28665.96754147422; S=0.3641122827668924; T=0.3641122827668925 }
    6  = @{ H=70.73212899999999; I=0.5111277390235027; J=0.5111277390235027;
             O=0.2125756143240238; P=0.2125756143240238;
             Q=950.4552186215608; R=8554.096967594045; S=0.1086532931209704; T=0.1086532931209704 }
    7  = @{ H=70.73212899999999; I=0.5111277390235027; J=0.5111277390235027;
             M=14.23299766666667; N=42.698993;         O=0.07505396441392481; P=0.07505396441392483;
             Q=335.5767423384552; R=3020.190681046096; S=0.03836216313563982; T=0.03836216313563983 }
    8  = @{ G=21.446869; H=64.34060699999999;           I=0.4649410310173153;  J=0.4649410310173154;
             M=135.0916853333333; N=405.2750559999999; O=0.7123704212620513;  P=0.7123704212620514;
             Q=2897.29367833322; R=26075.64310499898;  S=0.3312102381278174; T=0.3312102381278175 }
    9  = @{ G=21.446869; H=64.34060699999999;           I=0.4649410310173153;  J=0.4649410310173154;
             O=0.2125756143240238; P=0.2125756143240238;
             Q=864.5698433936425; R=7781.128590542782; S=0.0988351252929508; T=0.09883512529295083 }
    10 = @{ G=21.446869; H=64.34060699999999;           I=0.4649410310173153;  J=0.4649410310173154;
             M=14.23299766666667; N=42.698993;         O=0.07505396441392481; P=0.07505396441392483;
             Q=305.2532364343056; R=2747.279127908751; S=0.03489566759654709; T=0.03489566759654711 }
}

foreach ($rowNum in $updates.Keys) {
    $row = $updates[$rowNum]
    foreach ($col in $row.Keys) {
        $cellRef = "$col$rowNum"
        $ws.Range($cellRef).Value = $row[$col]
    }
}
